# Correção de URL de fotos
# The photo URL for NELSON LUIZ SPERLE TEICH (row 27, column B) pointed to a
# Wikipedia "Ficheiro:" page (not a direct image link) and is replaced with a
# working, direct image URL. The stale hyperlink that Excel had attached to
# the old (now wrong) address is removed as well, while the blue/underlined
# "Hiperlink" cell style is left in place (matches the target workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Replace the broken URL text in B27 with the corrected, direct image URL.
$ws.Range("B27").Value = "https://static.ndmais.com.br/2020/05/49811303336_8f58832171_c.jpg"

# 2. Remove the now-stale hyperlink relationship that still targeted the old
#    broken Wikipedia "Ficheiro:" URL (the hyperlinks on B16 and B2 stay).
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq "`$B`$27") {
        $hl.Delete()
    }
}

# 3. Mirror the author's view state: they scrolled down and the cell they
#    edited (B27) ended up selected/active.
$ws.Range("B27").Select()
